$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing holdings rows 2-13: Current Price, Profit/Loss,
#     Percentage Change, Todays Change, Todays Change Percent (columns F-J) ---
$ws.Cells.Item(2, 6).Value = 7.929999828338623
$ws.Cells.Item(2, 7).Value = -20101.00414498177
$ws.Cells.Item(2, 8).Value = -75.00787951989089
$ws.Cells.Item(2, 9).Value = 0.09999990463256836
$ws.Cells.Item(2, 10).Value = 1.28
$ws.Cells.Item(3, 6).Value = 41.34999847412109
$ws.Cells.Item(3, 7).Value = -18302.20090789795
$ws.Cells.Item(3, 8).Value = -42.65705384257233
$ws.Cells.Item(3, 9).Value = 1.669998168945312
$ws.Cells.Item(3, 10).Value = 4.21
$ws.Cells.Item(4, 6).Value = 70.73090362548828
$ws.Cells.Item(4, 7).Value = -425.0918655395502
$ws.Cells.Item(4, 8).Value = -3.320252015461612
$ws.Cells.Item(4, 9).Value = 0.4509048461914062
$ws.Cells.Item(4, 10).Value = 0.64
$ws.Cells.Item(5, 6).Value = 8.569999694824219
$ws.Cells.Item(5, 7).Value = 4873.399369812012
$ws.Cells.Item(5, 8).Value = 38.00321569765248
$ws.Cells.Item(5, 9).Value = -0.130000114440918
$ws.Cells.Item(5, 10).Value = -1.49
$ws.Cells.Item(6, 6).Value = 135.4100036621094
$ws.Cells.Item(6, 7).Value = -2654.578927001953
$ws.Cells.Item(6, 8).Value = -6.271195637772981
$ws.Cells.Item(6, 9).Value = 4.270004272460938
$ws.Cells.Item(6, 10).Value = 3.26
$ws.Cells.Item(7, 6).Value = 251.1300048828125
$ws.Cells.Item(7, 7).Value = 734.5307394140627
$ws.Cells.Item(7, 8).Value = 2.527151499474362
$ws.Cells.Item(7, 9).Value = 2.620010375976562
$ws.Cells.Item(7, 10).Value = 1.05
$ws.Cells.Item(8, 6).Value = 382.2099914550781
$ws.Cells.Item(8, 7).Value = 1335.525305090331
$ws.Cells.Item(8, 8).Value = 18.08996831708524
$ws.Cells.Item(8, 9).Value = 3.25
$ws.Cells.Item(8, 10).Value = 0.86
$ws.Cells.Item(9, 6).Value = 177.0099945068359
$ws.Cells.Item(9, 7).Value = -48.50013732910128
$ws.Cells.Item(9, 8).Value = -1.084104774050881
$ws.Cells.Item(9, 9).Value = 0.5099945068359375
$ws.Cells.Item(9, 10).Value = 0.29
$ws.Cells.Item(10, 6).Value = 79.33499908447266
$ws.Cells.Item(10, 7).Value = 0.4751238400268634
$ws.Cells.Item(10, 8).Value = 0.4749228526756111
$ws.Cells.Item(10, 9).Value = 0.08499908447265625
$ws.Cells.Item(10, 10).Value = 0.11
$ws.Cells.Item(11, 6).Value = 125.5192031860352
$ws.Cells.Item(11, 7).Value = 1.037808564758298
$ws.Cells.Item(11, 8).Value = 1.037755120369598
$ws.Cells.Item(11, 9).Value = 0.2392044067382812
$ws.Cells.Item(11, 10).Value = 0.19
$ws.Cells.Item(12, 6).Value = 28.77239990234375
$ws.Cells.Item(12, 7).Value = -0.5448235375976552
$ws.Cells.Item(12, 8).Value = -0.5447635591297951
$ws.Cells.Item(12, 9).Value = 0.2323989868164062
$ws.Cells.Item(12, 10).Value = 0.8100000000000001
$ws.Cells.Item(13, 6).Value = 48.03979873657227
$ws.Cells.Item(13, 7).Value = -1.838210979919429
$ws.Cells.Item(13, 8).Value = -1.839397759353764
$ws.Cells.Item(13, 9).Value = 0.3097991943359375
$ws.Cells.Item(13, 10).Value = 0.65

# --- Append two new holdings: AMZN (stock) and VOO (etf) ---
$ws.Cells.Item(14, 1).Value = "stock"
$ws.Cells.Item(14, 2).Value = "AMZN"
$ws.Cells.Item(14, 3).Value = 236.32
$ws.Cells.Item(14, 4).Value = 4
$ws.Cells.Item(14, 5).Value = 945.28
$ws.Cells.Item(14, 6).Value = 229.6900024414062
$ws.Cells.Item(14, 7).Value = -26.51999023437497
$ws.Cells.Item(14, 8).Value = -2.805516908680494
$ws.Cells.Item(14, 9).Value = 0.760009765625
$ws.Cells.Item(14, 10).Value = 0.33

$ws.Cells.Item(15, 1).Value = "etf"
$ws.Cells.Item(15, 2).Value = "VOO"
$ws.Cells.Item(15, 3).Value = 556.53
$ws.Cells.Item(15, 4).Value = 0.44921
$ws.Cells.Item(15, 5).Value = 249.9988413
$ws.Cells.Item(15, 6).Value = 558.5549926757812
$ws.Cells.Item(15, 7).Value = 0.9096469598877076
$ws.Cells.Item(15, 8).Value = 0.3638604703755911
$ws.Cells.Item(15, 9).Value = 3.7449951171875
$ws.Cells.Item(15, 10).Value = 0.68
